$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the description text for "num_acessos" (row 15, column B) from
# "número de clientes " to "número de clientes banda larga fixa total"
$ws.Range("B15").Value = "número de clientes banda larga fixa total"

# Move the active selection to B16, matching the saved view state
$ws.Range("B16").Select()
